{"js": "const body = context.document.body;\n\n// 1) \" for a Trial to the Court.\" -> \" for a trial to the Court.\"\n//    (de-capitalize \"Trial\")\nlet results = body.search(\" for a Trial to the Court.\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\" for a trial to the Court.\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Remove the bold \"Offense of Violence. \" lead-in and instead say\n//    \"The Court additionally informed the Defendant \".\nresults = body.search(\"Offense of Violence. \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"\", \"Replace\");\n  await context.sync();\n}\n\nresults = body.search(\"The Court informed the Defendant \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"The Court additionally informed the Defendant \", \"Replace\");\n  await context.sync();\n}\n\n// 3) \"{{'\\n'}}The Court heard statements...\" -> \"{{'\\n'}}Prior to sentencing, the Court heard statements...\"\nresults = body.search(\"The Court heard statements from the Prosecutor\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Prior to sentencing, the Court heard statements from the Prosecutor\", \"Replace\");\n  await context.sync();\n}\n\n// 4) \" fees monthly.  \" -> \" fees monthly.\" (drop the trailing double space)\nresults = body.search(\"fees monthly.  \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"fees monthly.\", \"Replace\");\n  await context.sync();\n}\n\n// 5) Remove the stray trailing space run after \"... immediately.\"\nresults = body.search(\"immediately. \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"immediately.\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Once($searchText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Wrap = 0\n    $found = $find.Execute()\n    if ($found) {\n        $rng.Text = $replaceText\n    } else {\n        Write-Output \"NOT FOUND: $searchText\"\n    }\n    return $found\n}\n\n# 1) \" for a Trial to the Court.\" -> \" for a trial to the Court.\"\n#    (de-capitalize \"Trial\")\nReplace-Once \" for a Trial to the Court.\" \" for a trial to the Court.\"\n\n# 2) Remove the bold \"Offense of Violence. \" lead-in and instead say\n#    \"The Court additionally informed the Defendant \".\nReplace-Once \"Offense of Violence. \" \"\"\nReplace-Once \"The Court informed the Defendant \" \"The Court additionally informed the Defendant \"\n\n# 3) \"{{'\\n'}}The Court heard statements...\" -> \"{{'\\n'}}Prior to sentencing, the Court heard statements...\"\nReplace-Once \"The Court heard statements from the Prosecutor\" \"Prior to sentencing, the Court heard statements from the Prosecutor\"\n\n# 4) \" fees monthly.  \" -> \" fees monthly.\" (drop the trailing double space)\nReplace-Once \" fees monthly.  \" \" fees monthly.\"\n\n# 5) Remove the stray trailing space run after \"... immediately.\"\nReplace-Once \"immediately. \" \"immediately.\"\n"}
